# tasks_db.xlsx: turn the three " - Task N" rows (each dated 2025-01-20, with
# a custom "YYYY-MM-DD HH:MM:SS" number format) into a plain two-row table
# with a real header row ("Date" / "Tasks") and a single sample data row
# ("2025-01-31" / "test").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 1: header labels -------------------------------------------------
# B1 already carries the bold / bordered / centered header style - keep it.
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Tasks"

# A1 used to be styled as a date cell; restyle it to match B1's header look
# by copying B1's formatting over (format-only paste, so the new text stays).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial($xlPasteFormats)

# --- Row 2: single plain sample row ---------------------------------------
# Enter "2025-01-31" as literal text (not an auto-converted date serial) by
# switching the cell to the Text number format before typing it in.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-01-31"

# Now drop A2 back to an unstyled (General) look by copying B2's current
# (default) formatting over it, same trick as above.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)

$ws.Range("B2").Value = "test"

# --- Drop the old third task row entirely ---------------------------------
$ws.Rows.Item(3).Delete()

$excel.CutCopyMode = $false
